$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.859.65'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.525.52'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.17'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '197.53'
$ws.Range('E6').Value = '  +6.48%  '
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.198'
$ws.Range('E9').Value = '  -7.94%  '
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.82'
$ws.Range('E11').Value = '  +1.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000302'
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.085.02'
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '598.99'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '70.022.85'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.09'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.75'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.525.25'
$ws.Range('E19').Value = '  +0.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.122'
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.993'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.36'
$ws.Range('E22').Value = '  +7.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.29'
$ws.Range('E23').Value = '  +4.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '101.97'
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.63'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.21'
$ws.Range('E26').Value = '  +6.51%  '
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.62'
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.36'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('E30').Value = '  +11.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.07'
$ws.Range('E31').Value = '  +1.61%  '
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.16'
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0863'
$ws.Range('E35').Value = '  +11.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.724.28'
$ws.Range('E36').Value = '  +3.47%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.07'
$ws.Range('E38').Value = '  -3.69%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.393'
$ws.Range('E40').Value = '  -1.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.64'
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '488.53'
$ws.Range('E42').Value = '  -5.80%  '
$ws.Range('E43').Value = '  -3.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0454'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('E45').Value = '  -3.33%  '
$ws.Range('E46').Value = '  -2.00%  '
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('E49').Value = '  -3.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000250'
$ws.Range('E50').Value = '  +2.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '129.71'
$ws.Range('E51').Value = '  -0.90%  '
